# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.150.54"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "2.468.55"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.18"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.15"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.10%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.512"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("E9").Value = "  +2.41%  "

$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.94"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").Value = "2.928.75"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.41"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").Value = "67.151.50"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "2.383.91"
$ws.Range("E17").Value = "  -3.76%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.90"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -2.73%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -1.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.03"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.19"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.19"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.79"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.12"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").Value = "2.597.37"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").Value = "0.0₃0901"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "500.85"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.73"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.01"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.75%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.32"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  +1.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.327"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.82"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "141.98"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.511"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").Value = "0.0₆0254"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.57"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -1.58%  "

$ws.Range("E51").Value = "  +0.02%  "

